$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 (second data row), pushing
# the existing row 3 data down to row 4. Excel's default row-insert picks up
# the formatting of the row above (row 2), which is why the new row 3 keeps
# the style used for the hyperlink cell E2 (style index 2) on cell E3.
$null = $ws.Rows("3:3").Insert()

# Leave the selection on row 2, matching the post-edit cursor position.
$null = $ws.Range("A1").Select()
$null = $ws.Rows("2:2").Select()
